$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Fix incorrect logical value export from CSV: several fields in the data
# dictionary were mis-typed as "character" in column I ("type"); correct
# them to the type that actually matches each ODK field.
$ws.Range("I2").Value  = "date"       # date
$ws.Range("I3").Value  = "timestamp"  # start
$ws.Range("I4").Value  = "timestamp"  # end
$ws.Range("I12").Value = "date"       # a1-enroldate
$ws.Range("I39").Value = "date"       # o1-o1_2
$ws.Range("I56").Value = "date"       # n1-n1_6
$ws.Range("I58").Value = "integer"    # n1-maxduration
$ws.Range("I60").Value = "integer"    # n1-n1_8

# Column I needs to be a touch wider to fit the longer new values.
$ws.Columns.Item(9).ColumnWidth = 9.43

# Saved view state: zoom out and scroll/select further down the sheet.
$win.Zoom = 50
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("N36").Select()
